# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Each entry below maps a row number to its new F-column value for a given sheet.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "展览" = @{
        2  = 151
        4  = 38
        5  = 151
        11 = 6
        13 = 167
        18 = 5113
        19 = 58
        22 = 2279
        25 = 2128
    }
    "全部类型" = @{
        2  = 151
        4  = 38
        5  = 151
        11 = 6
        13 = 167
        18 = 5113
        20 = 58
        24 = 2279
        27 = 2128
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}

$wb.Save()
